{"js": "const replacements = [\n  [\"530\u00f75=\", \"421\u00f78=\"],\n  [\"215\u00f73=\", \"847\u00f76=\"],\n  [\"788\u00f76=\", \"535\u00f78=\"],\n  [\"380\u00f77=\", \"776\u00f73=\"],\n  [\"777\u00f77=\", \"988\u00f78=\"],\n  [\"122\u00f73=\", \"977\u00f76=\"],\n  [\"683\u00f77=\", \"160\u00f78=\"],\n  [\"271\u00f77=\", \"638\u00f74=\"],\n  [\"486\u00f78=\", \"104\u00f78=\"],\n  [\"329\u00f79=\", \"279\u00f78=\"],\n  [\"361\u00f79=\", \"684\u00f74=\"],\n  [\"162\u00f77=\", \"559\u00f76=\"],\n  [\"363\u00f76=\", \"834\u00f73=\"],\n  [\"654\u00f73=\", \"445\u00f72=\"],\n  [\"256\u00f75=\", \"991\u00f79=\"],\n  [\"852\u00f78=\", \"609\u00f75=\"],\n  [\"527\u00f75=\", \"610\u00f73=\"],\n  [\"243\u00f74=\", \"252\u00f75=\"],\n  [\"445\u00f79=\", \"598\u00f77=\"],\n  [\"669\u00f75=\", \"180\u00f74=\"],\n  [\"267\u00f79=\", \"517\u00f79=\"],\n  [\"801\u00f79=\", \"957\u00f73=\"],\n  [\"249\u00f77=\", \"347\u00f74=\"],\n  [\"776\u00f78=\", \"483\u00f72=\"],\n  [\"121\u00f74=\", \"742\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"530\u00f75=\", \"421\u00f78=\"),\n    @(\"215\u00f73=\", \"847\u00f76=\"),\n    @(\"788\u00f76=\", \"535\u00f78=\"),\n    @(\"380\u00f77=\", \"776\u00f73=\"),\n    @(\"777\u00f77=\", \"988\u00f78=\"),\n    @(\"122\u00f73=\", \"977\u00f76=\"),\n    @(\"683\u00f77=\", \"160\u00f78=\"),\n    @(\"271\u00f77=\", \"638\u00f74=\"),\n    @(\"486\u00f78=\", \"104\u00f78=\"),\n    @(\"329\u00f79=\", \"279\u00f78=\"),\n    @(\"361\u00f79=\", \"684\u00f74=\"),\n    @(\"162\u00f77=\", \"559\u00f76=\"),\n    @(\"363\u00f76=\", \"834\u00f73=\"),\n    @(\"654\u00f73=\", \"445\u00f72=\"),\n    @(\"256\u00f75=\", \"991\u00f79=\"),\n    @(\"852\u00f78=\", \"609\u00f75=\"),\n    @(\"527\u00f75=\", \"610\u00f73=\"),\n    @(\"243\u00f74=\", \"252\u00f75=\"),\n    @(\"445\u00f79=\", \"598\u00f77=\"),\n    @(\"669\u00f75=\", \"180\u00f74=\"),\n    @(\"267\u00f79=\", \"517\u00f79=\"),\n    @(\"801\u00f79=\", \"957\u00f73=\"),\n    @(\"249\u00f77=\", \"347\u00f74=\"),\n    @(\"776\u00f78=\", \"483\u00f72=\"),\n    @(\"121\u00f74=\", \"742\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2, [ref]$false, [ref]$false, [ref]$false, [ref]$false) | Out-Null\n}\n"}
